$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows being appended to the "location" list (rows 233-235),
# matching columns: A=Category, B=latitude/longitude, C=Location,
# D=CITY, E=COUNTRY, F=YouTube Link
#
# Cell values are written in the same order the new entries were typed
# in originally (Location, then YouTube link, then coordinates, then
# category), so newly-introduced shared strings land at the same index
# positions as in the authored workbook.

# --- Row 233 ---
$ws.Range("C233").Value = "【東京湾ライブカメラ】 LIVE Tokyo Bay, Traffic between ships and planes"
$ws.Range("F233").Value = "qMDxy_qbdtE"
$ws.Range("B233").Value = "35.630819822372594, 139.7704128344258"
$ws.Range("A233").Value = "LIVE, SEA, HARBOR"
$ws.Range("D233").Value = "Tokyo"
$ws.Range("E233").Value = "Japan"

# --- Row 234 ---
$ws.Range("C234").Value = "【有明ライブカメラ】 LIVE Scenery overlooking land"
$ws.Range("F234").Value = "_ByNEL0Ton4"
$ws.Range("B234").Value = "35.63151956157848, 139.78222536039993"
$ws.Range("A234").Value = "LIVE, SEA, HARBOR, BRIDGE"
$ws.Range("D234").Value = "Tokyo"
$ws.Range("E234").Value = "Japan"

# --- Row 235 ---
$ws.Range("C235").Value = "【お台場ライブカメラ】 Live Beach Cam - Obaiba Beach"
$ws.Range("B235").Value = "35.63192580012683, 139.7764094005543"
$ws.Range("F235").Value = "KsoxRtx01KE"
$ws.Range("A235").Value = "LIVE, SEA, BEACH, BUILDING"
$ws.Range("D235").Value = "Tokyo"
$ws.Range("E235").Value = "Japan"

# The rest of the list uses a thin left/right border on columns A and E
# (no fill). Re-apply that existing formatting to the new rows by
# copying it down from the row directly above, rather than rebuilding
# it from scratch, so no new/duplicate style entries get created.
$ws.Range("A232").Copy()
$ws.Range("A233:A235").PasteSpecial(-4122)
$ws.Range("E232").Copy()
$ws.Range("E233:E235").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Keep frozen pane / selection consistent with appended rows (Excel
# typically scrolls the view so the freshly-added row is visible and
# selects the next empty row underneath the new data).
$ws.Application.ActiveWindow.ScrollRow = 214
$ws.Range("A236").Select()
